$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.256.15"
$ws.Range("D3").Value = "1.594.18"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'211.61"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "'18.93"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "1.818.24"
$ws.Range("D13").Value = "1.605.75"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "'63.58"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "26.217.26"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "'229.23"
$ws.Range("E18").Value = "  +7.21%  "
$ws.Range("D19").Value = "0.0₃0721"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'7.63"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "'145.52"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'6.98"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'15.31"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "1.464.59"
$ws.Range("E33").Value = "  +4.02%  "
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "'0.564"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "'0.820"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "'5.78"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "'0.930"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "1.730.94"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "'0.754"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'87.53"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E51").Value = "  -2.47%  "
